$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-10 Saturday" "2024-08-11 Sunday"

Replace-Text "753×3=2259" "288×8=2304"
Replace-Text "347×4=1388" "449×9=4041"
Replace-Text "723×8=5784" "675×5=3375"
Replace-Text "588×8=4704" "433×6=2598"
Replace-Text "833×9=7497" "541×4=2164"

Replace-Text "160×2=320" "250×8=2000"
Replace-Text "250×9=2250" "883×2=1766"
Replace-Text "866×9=7794" "211×2=422"
Replace-Text "630×8=5040" "746×6=4476"
Replace-Text "881×4=3524" "306×2=612"

Replace-Text "512×9=4608" "385×6=2310"
Replace-Text "108×6=648" "163×8=1304"
Replace-Text "709×8=5672" "755×8=6040"
Replace-Text "427×2=854" "538×7=3766"
Replace-Text "221×6=1326" "255×7=1785"

Replace-Text "299×8=2392" "625×8=5000"
Replace-Text "701×3=2103" "587×9=5283"
Replace-Text "260×2=520" "817×6=4902"
Replace-Text "546×4=2184" "292×7=2044"
Replace-Text "508×8=4064" "679×5=3395"

Replace-Text "212×6=1272" "982×4=3928"
Replace-Text "617×3=1851" "692×8=5536"
Replace-Text "538×3=1614" "881×3=2643"
Replace-Text "740×2=1480" "638×9=5742"
Replace-Text "805×7=5635" "504×2=1008"
